$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.066983461380005
$ws.Range("B1").Value = 1.425067543983459
$ws.Range("C1").Value = 2.318607807159424
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.798639416694641
